# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 142 (pushing existing rows down by one)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 142; all data below shifts down.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record's values.
$ws.Cells.Item(142, 1).Value = 11
$ws.Cells.Item(142, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(142, 3).Value = "Bíobío"
$ws.Cells.Item(142, 4).NumberFormat = $ws.Cells.Item(143, 4).NumberFormat
$ws.Cells.Item(142, 4).Value = 45007
$ws.Cells.Item(142, 5).Value = 8
$ws.Cells.Item(142, 6).Value = "Fruta"
$ws.Cells.Item(142, 7).Value = 100109
$ws.Cells.Item(142, 8).Value = "Uva"
$ws.Cells.Item(142, 9).Value = 100109001
$ws.Cells.Item(142, 10).Value = "Uva"
$ws.Cells.Item(142, 11).Value = "Red Globe"
$ws.Cells.Item(142, 12).Value = "Primera"
$ws.Cells.Item(142, 13).Value = 180
$ws.Cells.Item(142, 14).Value = 10000
$ws.Cells.Item(142, 15).Value = 11000
$ws.Cells.Item(142, 16).Value = 10444
$ws.Cells.Item(142, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(142, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(142, 19).Value = 580
$ws.Cells.Item(142, 20).Value = 18
